# Fruta / hortaliza, semanal
#
# The weekly data refresh reshuffles the per-record columns (Fecha, Calidad,
# Volumen, Precio minimo/maximo/promedio, Unidad de comercializacion,
# Precio $/Kg, Kg/unidad) across the existing data rows (rows 2-26) of the
# single worksheet. The "identity" columns (Mercado, Region, Codreg, Tipo,
# Producto, Categoria, Variedad, Origen -> columns A,B,C,E,F,G,H,I,J,K,R)
# stay the same for every row, only the 9 "variable" columns move between
# rows according to a fixed permutation.
#
# Because several rows trade values with each other, we must snapshot all
# of the current values first and only then write the new values back -
# otherwise we would overwrite data before it has been read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers (1-based) of the fields that get reshuffled between rows.
# D=4 Fecha, L=12 Calidad, M=13 Volumen, N=14 Precio minimo,
# O=15 Precio maximo, P=16 Precio promedio ponderado,
# Q=17 Unidad de comercializacion, S=19 Precio $/Kg, T=20 Kg / unidad
$cols = @(4, 12, 13, 14, 15, 16, 17, 19, 20)

$firstRow = 2
$lastRow = 26

# 1) Snapshot every current value for each data row / reshuffled column.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# 2) Target row -> source row mapping describing where each row's new
#    values come from (derived from the committed workbook contents).
$mapping = @{
    2 = 22
    3 = 13
    4 = 14
    5 = 21
    6 = 24
    7 = 8
    8 = 9
    9 = 5
    10 = 6
    11 = 10
    12 = 15
    13 = 16
    14 = 11
    15 = 25
    16 = 12
    17 = 23
    18 = 2
    19 = 3
    20 = 26
    21 = 18
    22 = 7
    23 = 17
    24 = 20
    25 = 19
    26 = 4
}

# 3) Write the reshuffled values back using the snapshot so rows that swap
#    values with each other end up correct regardless of iteration order.
foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $src = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value = $src[$c]
    }
}
